$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (kNN) updates
$ws.Range("B2").Value = 0.6650652024708303
$ws.Range("C2").Value = 0.727838722474945
$ws.Range("D2").Value = 0.6650652024708303
$ws.Range("E2").Value = 0.6637261923393825
$ws.Range("N2").Value = 0.7874857012125374
$ws.Range("O2").Value = 0.800592546381022
$ws.Range("P2").Value = 0.7874857012125374
$ws.Range("Q2").Value = 0.7853567617702084

# Row 6 (Ensemble) updates
$ws.Range("G6").Value = 0.8789981776696546
$ws.Range("I6").Value = 0.8734765654573513
$ws.Range("J6").Value = 0.8283687943262411
$ws.Range("K6").Value = 0.8429374553699459
$ws.Range("L6").Value = 0.8283687943262411
$ws.Range("M6").Value = 0.8306296978538832
$ws.Range("O6").Value = 0.8629328911997162
$ws.Range("Q6").Value = 0.8577187713252601
$ws.Range("V6").Value = 0.8518874399450928
$ws.Range("W6").Value = 0.8559014625503328
$ws.Range("X6").Value = 0.8518874399450928
$ws.Range("Y6").Value = 0.8512286782416443
